$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style from B1 (bold font, thin border, centered/top aligned) to the
# two new header cells so they match the rest of row 1 (style index 1)
$ws.Range("B1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

function ConvertTo-Excel2DArray {
    param([string]$csv)
    $lines = $csv -split "`n" | Where-Object { $_.Length -gt 0 }
    $rowCount = $lines.Count
    $colCount = ($lines[0] -split ",").Count
    $arr = New-Object 'object[,]' $rowCount, $colCount
    for ($r = 0; $r -lt $rowCount; $r++) {
        $fields = $lines[$r] -split ","
        for ($c = 0; $c -lt $colCount; $c++) {
            $arr[$r, $c] = [double]$fields[$c]
        }
    }
    return $arr
}

# Recomputed results for columns B,C,D,E,F,G,H,I,J,K (rows 2:25)
$bkCsv = @"
2.817632658248954,0.9938639734979517,0,1.370168940525417,1.45952168114772,0.8269549254233368,0.01297449529001643,0.0080534191016155,0.5122514496681134,0.3547313468008468
2.449499810661564,0.8660181943978387,0,1.192060143739269,1.26082040731464,0.7404094690176493,0.009211099144520327,0.005422390222803219,0.4762333812276012,0.3379904458039533
2.222989174503198,0.7880524250603855,0,1.083167400509197,1.139972825279813,0.6878666164700888,0.007159029862464916,0.004065610465914027,0.4544845323456173,0.3279706305622447
2.128776028006882,0.7569848775538617,0,1.038855796295891,1.090973950927975,0.6654150457429893,0.006379783292803309,0.003624448366656541,0.4451103320581495,0.3233116060869605
2.110961852465437,0.7525564764336536,0,1.031469244301022,1.082851457025441,0.6602678000264888,0.006250811725672656,0.003625885627272396,0.4428336773037813,0.3217677630312821
2.215792765885737,0.7896158984223689,0,1.082476937428638,1.139311062524612,0.6836582212733902,0.007139509106109831,0.004238045131493173,0.4523753652420339,0.3257914928644965
2.682841354346067,0.9523075504260134,0,1.308513878348421,1.390742831773508,0.7917715554769984,0.01160558447156113,0.007293233014741496,0.4971154246358225,0.3460909881081236
3.606894177290826,1.27298481908872,0,1.756962358112077,1.895061171842798,1.016707132861157,0.02270243446163134,0.01535255824824855,0.5918489414939927,0.3921364426070966
4.255175590793669,1.497532504463607,0,1.987594654288458,2.275587181990815,1.160757846438145,0.03176002381723375,0.02284837784439908,0.6513058895570509,0.4180526566706817
4.298837547128755,1.460422464860471,0,1.291738902706456,2.451628164639231,1.015709253472465,0.04657845409392891,0.02430516621366774,0.5747520763726897,0.353366902713006
4.208316715658952,1.381149853074589,0,0.7897693304516196,2.518780942186339,0.877957098950148,0.08184144539288951,0.02378135406691051,0.5045464823185455,0.3008504647490682
4.0010314481097,1.268935617212605,0,0.4082876288041959,2.504295469238414,0.7324168996077702,0.1341392572930857,0.02198983365296669,0.4323116283865147,0.2514094177072579
3.802973308390278,1.177365325101846,0,0.2145324879897004,2.457142640142251,0.6283682704635964,0.181195998049958,0.02031809232954451,0.3813722449953758,0.2185953024860083
3.728028685284357,1.148631601191539,0,0.1762909600798928,2.428326063136538,0.599371274376935,0.1929418949175954,0.01972159508092464,0.3675878794254714,0.2102451943694348
3.494633174730097,1.080021390926504,0,0.1711578881997511,2.264147245220144,0.5691221847812642,0.1774795499779316,0.01710536318870393,0.3573465955604718,0.2090345593314673
3.41957652915039,1.075268760700453,0,0.2560380171113152,2.1642275708721,0.5995633359360824,0.1387227049625466,0.01597732866666135,0.375872776809075,0.2249295224733459
3.477429043483767,1.124922562326105,0,0.4918037635947954,2.107028586650728,0.6928362916651594,0.08678053725635948,0.01577058936645237,0.4245761180827685,0.2607544666684944
3.629449803906937,1.221748998862324,0,0.9308407959301093,2.087706772602971,0.8322995477357011,0.04474227722054991,0.01689155908493856,0.4947129355708455,0.3107020924110415
4.064611101601997,1.443907054595627,0,1.922111612942899,2.174835551210137,1.109194401771504,0.02915449852148821,0.02126713301065486,0.6287714278265355,0.4041376138654584
4.609475764742001,1.642715345676038,0,2.261661826454414,2.470978708896922,1.262348198430175,0.03818790976198017,0.02806413005870034,0.6954230391121143,0.4399794705342259
4.960282458136078,1.762862427750576,0,2.432581515462218,2.667416158031983,1.356805558019829,0.04412220771101616,0.03261362290041525,0.736309118782799,0.4622499029566356
4.780313082448458,1.696179896713318,0,2.341343363195989,2.562285531758789,1.310961088819454,0.0409307067450202,0.02994206671436306,0.7167827944325893,0.4528252863981947
4.088752933374849,1.4501965674541,0,1.998793301716304,2.170038921572868,1.133195923795029,0.02974835399357456,0.02105681232821155,0.6407826228801241,0.4143962500007987
3.34599588185705,1.18918344614508,0,1.634722191659506,1.757113814925091,0.9478954920941334,0.01938917328357714,0.01317554641741658,0.5621617113341699,0.375518792896969
"@
$ws.Range("B2:K25").Value = (ConvertTo-Excel2DArray $bkCsv)

# Column M (rows 2:25) is now all zero
$mCsv = @"
0
0
0
0
0
0
0
0
0
0
0
0
0
0
0
0
0
0
0
0
0
0
0
0
"@
$ws.Range("M2:M25").Value = (ConvertTo-Excel2DArray $mCsv)

# New columns P (header 14) and Q (header 15), rows 1:25
$pqCsv = @"
14,15
0.6537529194058749,0
0.6720768277829308,0
0.6839657274974549,0
0.6892666504468394,0
0.6905286117217599,0
0.6850711312880975,0
0.6613849947866655,0
0.6182918653243448,0
0.5971470783592068,0
0.6539086658632485,0
0.7149600064223449,0
0.7829321917071184,0
0.8347912042425207,0
0.8482417799864237,0
0.8397883818614105,0
0.8058583353104751,0
0.7488076293816732,0
0.6877419410918577,0
0.6065903967960296,0
0.5813095129789474,0
0.5671178770490215,0
0.5731068088486424,0
0.5997948377826674,0
0.6314649991632209,0
"@
$ws.Range("P1:Q25").Value = (ConvertTo-Excel2DArray $pqCsv)

Write-Host "done"
